# Generate Report for Handback
# Refresh the "Latest Handback DateTime" column (column G) on each
# per-language sheet to the timestamp of this handback run.

$wb = $excel.ActiveWorkbook

$updates = @(
    @{ Sheet = "zh-cn"; Old = "2016-02-22 08:50:29"; New = "2016-02-22 08:58:35" },
    @{ Sheet = "de-de"; Old = "2016-02-22 08:50:40"; New = "2016-02-22 08:58:45" },
    @{ Sheet = "ja-jp"; Old = "2016-02-22 08:50:50"; New = "2016-02-22 08:58:57" },
    @{ Sheet = "zh-tw"; Old = "2016-02-22 08:51:01"; New = "2016-02-22 08:59:06" }
)

foreach ($u in $updates) {
    $ws = $wb.Worksheets.Item($u.Sheet)
    $usedRange = $ws.UsedRange
    $lastRow = $usedRange.Rows.Count
    $colRange = $ws.Range("G1:G" + $lastRow)
    $colRange.Replace($u.Old, $u.New) | Out-Null
}
